$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 21.50357910362046
$ws.Cells.Item(2, 3).Value = 8.498476304270989
$ws.Cells.Item(2, 4).Value = 7.375701638816752
$ws.Cells.Item(2, 6).Value = 41.42118514932706
$ws.Cells.Item(2, 7).Value = 49.39068229033535
$ws.Cells.Item(2, 8).Value = 19.31818353026877
$ws.Cells.Item(2, 10).Value = 10.45979591541273
$ws.Cells.Item(2, 12).Value = 11.89794649146415

$ws.Cells.Item(3, 2).Value = 21.07558978756189
$ws.Cells.Item(3, 3).Value = 8.141008748767657
$ws.Cells.Item(3, 4).Value = 7.373451566904109
$ws.Cells.Item(3, 6).Value = 41.47631814205839
$ws.Cells.Item(3, 7).Value = 49.37980865967608
$ws.Cells.Item(3, 8).Value = 19.37319980132258
$ws.Cells.Item(3, 10).Value = 10.48781973200847
$ws.Cells.Item(3, 12).Value = 11.88644142342564

$ws.Cells.Item(4, 2).Value = 20.81268762813189
$ws.Cells.Item(4, 3).Value = 7.911842994097497
$ws.Cells.Item(4, 4).Value = 7.372420464184991
$ws.Cells.Item(4, 6).Value = 41.52183358724174
$ws.Cells.Item(4, 7).Value = 49.39054011884131
$ws.Cells.Item(4, 8).Value = 19.41131466835402
$ws.Cells.Item(4, 10).Value = 10.5059610127194
$ws.Cells.Item(4, 12).Value = 11.88103260462438

$ws.Cells.Item(5, 2).Value = 20.70567337080461
$ws.Cells.Item(5, 3).Value = 7.816102831235017
$ws.Cells.Item(5, 4).Value = 7.372088892090403
$ws.Cells.Item(5, 6).Value = 41.54330301736111
$ws.Cells.Item(5, 7).Value = 49.39927563913751
$ws.Cells.Item(5, 8).Value = 19.42793310623301
$ws.Cells.Item(5, 10).Value = 10.51358936621194
$ws.Cells.Item(5, 12).Value = 11.87924668898688

$ws.Cells.Item(6, 2).Value = 20.68791543694826
$ws.Cells.Item(6, 3).Value = 7.800065709398393
$ws.Cells.Item(6, 4).Value = 7.372039200812389
$ws.Cells.Item(6, 6).Value = 41.54704409281781
$ws.Cells.Item(6, 7).Value = 49.40098912447255
$ws.Cells.Item(6, 8).Value = 19.43075808645267
$ws.Cells.Item(6, 10).Value = 10.51487029893845
$ws.Cells.Item(6, 12).Value = 11.87897545253141

$ws.Cells.Item(7, 2).Value = 20.81124370894981
$ws.Cells.Item(7, 3).Value = 7.910561221746167
$ws.Cells.Item(7, 4).Value = 7.372415633051451
$ws.Cells.Item(7, 6).Value = 41.52211131856968
$ws.Cells.Item(7, 7).Value = 49.39064028906727
$ws.Cells.Item(7, 8).Value = 19.41153439696541
$ws.Cells.Item(7, 10).Value = 10.50606293644946
$ws.Cells.Item(7, 12).Value = 11.88100682331229

$ws.Cells.Item(8, 2).Value = 21.35612089860948
$ws.Cells.Item(8, 3).Value = 8.377275242172436
$ws.Cells.Item(8, 4).Value = 7.374853370144585
$ws.Cells.Item(8, 6).Value = 41.43776779261648
$ws.Cells.Item(8, 7).Value = 49.38331224500249
$ws.Cells.Item(8, 8).Value = 19.33625131946078
$ws.Cells.Item(8, 10).Value = 10.46926494618317
$ws.Cells.Item(8, 12).Value = 11.89363721539454

$ws.Cells.Item(9, 2).Value = 22.41733141075207
$ws.Cells.Item(9, 3).Value = 9.21272114974281
$ws.Cells.Item(9, 4).Value = 7.382393586356109
$ws.Cells.Item(9, 6).Value = 41.36536489230501
$ws.Cells.Item(9, 7).Value = 49.50754843626001
$ws.Cells.Item(9, 8).Value = 19.22317380683806
$ws.Cells.Item(9, 10).Value = 10.40449029902848
$ws.Cells.Item(9, 12).Value = 11.93144661126954

$ws.Cells.Item(10, 2).Value = 23.18435478339092
$ws.Cells.Item(10, 3).Value = 9.774596149766891
$ws.Cells.Item(10, 4).Value = 7.389588650801089
$ws.Cells.Item(10, 6).Value = 41.36938761489877
$ws.Cells.Item(10, 7).Value = 49.68370635422657
$ws.Cells.Item(10, 8).Value = 19.16136949354679
$ws.Cells.Item(10, 10).Value = 10.36136230130177
$ws.Cells.Item(10, 12).Value = 11.96704117609955

$ws.Cells.Item(11, 2).Value = 23.52892958581267
$ws.Cells.Item(11, 3).Value = 10.01840748863353
$ws.Cells.Item(11, 4).Value = 7.39321544746197
$ws.Cells.Item(11, 6).Value = 41.38372206486154
$ws.Cells.Item(11, 7).Value = 49.7822667248894
$ws.Cells.Item(11, 8).Value = 19.13791423732228
$ws.Cells.Item(11, 10).Value = 10.34270245328982
$ws.Cells.Item(11, 12).Value = 11.98489872857164

$ws.Cells.Item(12, 2).Value = 23.65866040467784
$ws.Cells.Item(12, 3).Value = 10.10900273748801
$ws.Cells.Item(12, 4).Value = 7.394639161115293
$ws.Cells.Item(12, 6).Value = 41.39095204562795
$ws.Cells.Item(12, 7).Value = 49.82223025279507
$ws.Cells.Item(12, 8).Value = 19.12970555685449
$ws.Cells.Item(12, 10).Value = 10.33577374994888
$ws.Cells.Item(12, 12).Value = 11.99189728969107

$ws.Cells.Item(13, 2).Value = 23.6307559455052
$ws.Cells.Item(13, 3).Value = 10.08956891188661
$ws.Cells.Item(13, 4).Value = 7.394330309400566
$ws.Cells.Item(13, 6).Value = 41.38931476200855
$ws.Cells.Item(13, 7).Value = 49.81350612529995
$ws.Cells.Item(13, 8).Value = 19.13144344765189
$ws.Cells.Item(13, 10).Value = 10.33725986934713
$ws.Cells.Item(13, 12).Value = 11.99037956793554

$ws.Cells.Item(14, 2).Value = 23.5396183933096
$ws.Cells.Item(14, 3).Value = 10.02589568512481
$ws.Cells.Item(14, 4).Value = 7.393331571784659
$ws.Cells.Item(14, 6).Value = 41.38428075313877
$ws.Cells.Item(14, 7).Value = 49.78550168742987
$ws.Cells.Item(14, 8).Value = 19.13722539213323
$ws.Cells.Item(14, 10).Value = 10.34212967471501
$ws.Cells.Item(14, 12).Value = 11.98546978768127

$ws.Cells.Item(15, 2).Value = 23.48369241679649
$ws.Cells.Item(15, 3).Value = 9.986667637595122
$ws.Cells.Item(15, 4).Value = 7.392726353723527
$ws.Cells.Item(15, 6).Value = 41.38143201024192
$ws.Cells.Item(15, 7).Value = 49.76869174344645
$ws.Cells.Item(15, 8).Value = 19.14085477847008
$ws.Cells.Item(15, 10).Value = 10.34513044474255
$ws.Cells.Item(15, 12).Value = 11.98249307997379

$ws.Cells.Item(16, 2).Value = 23.16173749328734
$ws.Cells.Item(16, 3).Value = 9.758422103683289
$ws.Cells.Item(16, 4).Value = 7.389358702897969
$ws.Cells.Item(16, 6).Value = 41.36870280297158
$ws.Cells.Item(16, 7).Value = 49.67763532696532
$ws.Cells.Item(16, 8).Value = 19.16299639862488
$ws.Cells.Item(16, 10).Value = 10.36260102092148
$ws.Cells.Item(16, 12).Value = 11.96590738698931

$ws.Cells.Item(17, 2).Value = 22.96302224911496
$ws.Cells.Item(17, 3).Value = 9.61535421215571
$ws.Cells.Item(17, 4).Value = 7.387382973189148
$ws.Cells.Item(17, 6).Value = 41.3640997246087
$ws.Cells.Item(17, 7).Value = 49.62648957637532
$ws.Cells.Item(17, 8).Value = 19.17777542500946
$ws.Cells.Item(17, 10).Value = 10.37356394503437
$ws.Cells.Item(17, 12).Value = 11.95615694362694

$ws.Cells.Item(18, 2).Value = 22.84832436928253
$ws.Cells.Item(18, 3).Value = 9.531957472391536
$ws.Cells.Item(18, 4).Value = 7.386279912843633
$ws.Cells.Item(18, 6).Value = 41.36262902415845
$ws.Cells.Item(18, 7).Value = 49.59880697306212
$ws.Cells.Item(18, 8).Value = 19.18671435432103
$ws.Cells.Item(18, 10).Value = 10.37995985606823
$ws.Cells.Item(18, 12).Value = 11.95070570565116

$ws.Cells.Item(19, 2).Value = 22.80942456354814
$ws.Cells.Item(19, 3).Value = 9.503531679233062
$ws.Cells.Item(19, 4).Value = 7.385912175901703
$ws.Cells.Item(19, 6).Value = 41.36233305755718
$ws.Cells.Item(19, 7).Value = 49.58973231092683
$ws.Cells.Item(19, 8).Value = 19.18981612014025
$ws.Cells.Item(19, 10).Value = 10.38214093282897
$ws.Cells.Item(19, 12).Value = 11.94888706205603

$ws.Cells.Item(20, 2).Value = 22.98421838837214
$ws.Cells.Item(20, 3).Value = 9.630698971635224
$ws.Cells.Item(20, 4).Value = 7.387589847635819
$ws.Cells.Item(20, 6).Value = 41.36446789680458
$ws.Cells.Item(20, 7).Value = 49.63175460360269
$ws.Cells.Item(20, 8).Value = 19.17615677328181
$ws.Cells.Item(20, 10).Value = 10.37238757870718
$ws.Cells.Item(20, 12).Value = 11.95717867280377

$ws.Cells.Item(21, 2).Value = 23.56640908959456
$ws.Cells.Item(21, 3).Value = 10.0446452832824
$ws.Cells.Item(21, 4).Value = 7.393623563785304
$ws.Cells.Item(21, 6).Value = 41.38571044123375
$ws.Cells.Item(21, 7).Value = 49.79365568909402
$ws.Cells.Item(21, 8).Value = 19.13550879643562
$ws.Cells.Item(21, 10).Value = 10.34069557099253
$ws.Cells.Item(21, 12).Value = 11.98690552258528

$ws.Cells.Item(22, 2).Value = 23.94246660962553
$ws.Cells.Item(22, 3).Value = 10.30508218459134
$ws.Cells.Item(22, 4).Value = 7.397860019935139
$ws.Cells.Item(22, 6).Value = 41.41009665273576
$ws.Cells.Item(22, 7).Value = 49.91485313644519
$ws.Cells.Item(22, 8).Value = 19.11286895830432
$ws.Cells.Item(22, 10).Value = 10.32078348493419
$ws.Cells.Item(22, 12).Value = 12.00770942378324

$ws.Cells.Item(23, 2).Value = 23.74220365428754
$ws.Cells.Item(23, 3).Value = 10.1670165986759
$ws.Cells.Item(23, 4).Value = 7.395572307710896
$ws.Cells.Item(23, 6).Value = 41.39611946914986
$ws.Cells.Item(23, 7).Value = 49.84876400957733
$ws.Cells.Item(23, 8).Value = 19.12459202077644
$ws.Cells.Item(23, 10).Value = 10.33133788328474
$ws.Cells.Item(23, 12).Value = 11.99648121763823

$ws.Cells.Item(24, 2).Value = 22.9746370176768
$ws.Cells.Item(24, 3).Value = 9.623765169904971
$ws.Cells.Item(24, 4).Value = 7.387496217463379
$ws.Cells.Item(24, 6).Value = 41.36429778434017
$ws.Cells.Item(24, 7).Value = 49.62936892105022
$ws.Cells.Item(24, 8).Value = 19.17688718832066
$ws.Cells.Item(24, 10).Value = 10.37291912363403
$ws.Cells.Item(24, 12).Value = 11.95671626780775

$ws.Cells.Item(25, 2).Value = 22.13189551711581
$ws.Cells.Item(25, 3).Value = 8.995626346698158
$ws.Cells.Item(25, 4).Value = 7.38006108130599
$ws.Cells.Item(25, 6).Value = 41.37493575418834
$ws.Cells.Item(25, 7).Value = 49.45904360851442
$ws.Cells.Item(25, 8).Value = 19.25004291800042
$ws.Cells.Item(25, 10).Value = 10.42122706520004
$ws.Cells.Item(25, 12).Value = 11.91983508806199
